$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Inflammatory-Mac" cluster label to "Neutrophils" everywhere it's used
$ws.Cells.Replace("Inflammatory-Mac", "Neutrophils") | Out-Null

# Update row 2 values (Edge average/total expression derived specificity columns)
$ws.Range("I2").Value = 0.180007610261557
$ws.Range("J2").Value = 0.180007610261557
$ws.Range("S2").Value = 0.180007610261557
$ws.Range("T2").Value = 0.180007610261557

# Update row 3 values
$ws.Range("G3").Value = 0.1644956666666667
$ws.Range("H3").Value = 0.493487
$ws.Range("I3").Value = 0.819992389738443
$ws.Range("J3").Value = 0.819992389738443
$ws.Range("Q3").Value = 0.04583627886155556
$ws.Range("R3").Value = 0.4125265097540001
$ws.Range("S3").Value = 0.819992389738443
$ws.Range("T3").Value = 0.819992389738443

# Remove row 4 (the "Resolving-Mac" / Ccl22-Ackr2 pairing row) entirely
$ws.Rows("4:4").Delete() | Out-Null
